$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.076.12'
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").Value = '1.909.04'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("D4").Value = '''1.0000'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''333.28'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '''0.4643'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("D8").Value = '''0.4091'
$ws.Range("E8").Value = '  +3.26%  '
$ws.Range("D9").Value = '''47.89'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = '''0.08025'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").Value = '''1.008'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '''21.87'
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").Value = '1.909.46'
$ws.Range("E13").Value = '  +2.47%  '
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = '''7.103'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '''89.19'
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = '''0.9997'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").Value = '''0.06574'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").Value = '''17.55'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '''0.9991'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").Value = '29.084.58'
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").Value = '''5.453'
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = '''11.29'
$ws.Range("E24").Value = '  +2.22%  '
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").Value = '2.130.40'
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").Value = '''157.33'
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("D28").Value = '''19.77'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '''2.117'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '''5.427'
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").Value = '''119.09'
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").Value = '''0.9944'
$ws.Range("E32").Value = '  +2.55%  '
$ws.Range("D33").Value = '''0.09430'
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("D34").Value = '''1.428'
$ws.Range("E34").Value = '  +3.94%  '
$ws.Range("D35").Value = '''3.588'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = '''5.329'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = '''0.06111'
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = '''0.02247'
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("D39").Value = '''8.395'
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("D40").Value = '''1.178'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''0.5831'
$ws.Range("D42").Value = '''0.9990'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("E44").Value = '  -2.62%  '
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("D46").Value = '''2.317'
$ws.Range("E46").Value = '  +12.61%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5522'
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''12.10'
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.920'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.07044'
$ws.Range("E50").Value = '  +2.02%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '''47.61'
$ws.Range("E51").Value = '  +22.52%  '
